$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Add a new "2022-Q1" sheet, positioned right before "总计", by
#    copying the "2021-Q4" sheet (same column layout/styles) and then
#    overwriting its data with the 2022-Q1 figures.
# ---------------------------------------------------------------------
$srcSheet = $wb.Worksheets.Item("2021-Q4")
$totalSheet = $wb.Worksheets.Item("总计")
$srcSheet.Copy($totalSheet)

# The copy gets auto-named "2021-Q4 (2)" - grab it via its position
# (right before 总计) and rename it. Re-fetch 总计 since its Index
# shifted once the new sheet was inserted in front of it.
$totalSheet = $wb.Worksheets.Item("总计")
$newSheet = $wb.Worksheets.Item($totalSheet.Index - 1)
$newSheet.Name = "2022-Q1"

# Row 2
$newSheet.Range("B2").Value = "'005613"
$newSheet.Range("C2").Value = "上投摩根富时发达市场REITs指数QDII人民币份额"
$newSheet.Range("D2").Value = "'4.84"
$newSheet.Range("E2").Value = "'91.10"
$newSheet.Range("F2").Value = "'9.39"
$newSheet.Range("G2").Value = "'0.4545"
$newSheet.Range("H2").Value = 1

# Row 3
$newSheet.Range("B3").Value = "'005614"
$newSheet.Range("C3").Value = "上投摩根富时发达市场REITs指数QDII美钞"
$newSheet.Range("D3").Value = "'4.84"
$newSheet.Range("E3").Value = "'91.10"
$newSheet.Range("F3").Value = "'9.39"
$newSheet.Range("G3").Value = "'0.4545"
$newSheet.Range("H3").Value = 1

# Row 4
$newSheet.Range("B4").Value = "'005615"
$newSheet.Range("C4").Value = "上投摩根富时发达市场REITs指数QDII美汇"
$newSheet.Range("D4").Value = "'4.84"
$newSheet.Range("E4").Value = "'91.10"
$newSheet.Range("F4").Value = "'9.39"
$newSheet.Range("G4").Value = "'0.4545"
$newSheet.Range("H4").Value = 1

# Row 5
$newSheet.Range("B5").Value = "'000179"
$newSheet.Range("C5").Value = "广发美国房地产指数QDII-人民币"
$newSheet.Range("D5").Value = "'2.37"
$newSheet.Range("E5").Value = "'92.38"
$newSheet.Range("F5").Value = "'9.04"
$newSheet.Range("G5").Value = "'0.2142"
$newSheet.Range("H5").Value = 1

# Row 6
$newSheet.Range("B6").Value = "'000180"
$newSheet.Range("C6").Value = "广发美国房地产指数QDII - 美元"
$newSheet.Range("D6").Value = "'2.37"
$newSheet.Range("E6").Value = "'92.38"
$newSheet.Range("F6").Value = "'9.04"
$newSheet.Range("G6").Value = "'0.2142"
$newSheet.Range("H6").Value = 1

# Row 7
$newSheet.Range("B7").Value = "'160140"
$newSheet.Range("C7").Value = "南方道琼斯美国精选REIT指数(QDII-LOF)A"
$newSheet.Range("D7").Value = "'1.35"
$newSheet.Range("E7").Value = "'89.10"
$newSheet.Range("F7").Value = "'9.56"
$newSheet.Range("G7").Value = "'0.1291"
$newSheet.Range("H7").Value = 1

# Row 8
$newSheet.Range("B8").Value = "'006555"
$newSheet.Range("C8").Value = "浦银安盛全球智能科技股票（QDII）"
$newSheet.Range("D8").Value = "'3.20"
$newSheet.Range("E8").Value = "'85.41"
$newSheet.Range("F8").Value = "'2.77"
$newSheet.Range("G8").Value = "'0.0886"
$newSheet.Range("H8").Value = 5

# Row 9
$newSheet.Range("B9").Value = "'070031"
$newSheet.Range("C9").Value = "嘉实全球房地产(QDII)"
$newSheet.Range("D9").Value = "'0.60"
$newSheet.Range("E9").Value = "'95.08"
$newSheet.Range("F9").Value = "'7.71"
$newSheet.Range("G9").Value = "'0.0463"
$newSheet.Range("H9").Value = 1

# Row 10
$newSheet.Range("B10").Value = "'160141"
$newSheet.Range("C10").Value = "南方道琼斯美国精选REIT指数(QDII-LOF)C"
$newSheet.Range("D10").Value = "'0.44"
$newSheet.Range("E10").Value = "'89.10"
$newSheet.Range("F10").Value = "'9.56"
$newSheet.Range("G10").Value = "'0.0421"
$newSheet.Range("H10").Value = 1

# ---------------------------------------------------------------------
# 2) Update the "总计" sheet: insert a new row 2 for 2022-Q1, pushing
#    the earlier quarters down, and renumber the index column.
# ---------------------------------------------------------------------
$totalSheet.Rows.Item(2).Insert()

# The inserted row inherits formatting from the row above (the header);
# strip that from the data cells so they go back to the default style
# used by the other data rows.
$totalSheet.Range("B2:D2").ClearFormats()

# Give A2 the same style as the other index cells (bold/bordered).
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 9
$totalSheet.Range("D2").Value = 2.1

$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3
